# RPA datasets push 2024-06-14
# Updates the IPO tracking sheet:
#  - 에스오에스랩 (row 18) confirmed offer price (D18) is now set to 11500
#    (was "-", i.e. not yet determined) and the offering amount (E18)
#    is updated from 15000 to 23000.
#  - The lead underwriter list for 시프트업(유가) (row 21, column F) is
#    reformatted to use "." separators instead of "," between names.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# D18 stores the confirmed offer price as text (matching the "-" placeholder
# text that was there before), so force the cell to Text format first, then
# restore the original (default) cell style so no formatting change lingers.
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11500"
$ws.Range("D18").Style = "Normal"

$ws.Range("E18").Value = 23000
$ws.Range("F21").Value = "한국투자.NH투자.신한투자증권"
